# Generate Report for Handoff
# Rows 4-7 on the "zh-cn" and "de-de" sheets move from Priority "low" to "ht",
# and pick up a fresh "Latest Handoff Datetime" timestamp.

$wb = $excel.ActiveWorkbook

$ws_zhcn = $wb.Worksheets.Item("zh-cn")
for ($r = 4; $r -le 7; $r++) {
    $ws_zhcn.Cells.Item($r, 5).Value = "ht"
    $ws_zhcn.Cells.Item($r, 8).Value = "2016-11-09 02:07:13"
}

$ws_dede = $wb.Worksheets.Item("de-de")
for ($r = 4; $r -le 7; $r++) {
    $ws_dede.Cells.Item($r, 5).Value = "ht"
    $ws_dede.Cells.Item($r, 8).Value = "2016-11-09 02:07:29"
}
